$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 100, pushing the existing rows 100-129 down to 101-130.
$ws.Rows.Item(100).Insert()

# Populate the newly inserted row 100 with the new weekly record.
$ws.Range("A100").Value = 10
$ws.Range("B100").Value = "Vega Modelo de Temuco"
$ws.Range("C100").Value = "La Araucanía"
$ws.Range("D100").Value = 44463
$ws.Range("E100").Value = 9
$ws.Range("F100").Value = 100112005
$ws.Range("G100").Value = "Puerro"
$ws.Range("H100").Value = "Azul de Maquehue"
$ws.Range("I100").Value = "Primera"
$ws.Range("J100").Value = 30
$ws.Range("K100").Value = 8000
$ws.Range("L100").Value = 8000
$ws.Range("M100").Value = 8000
$ws.Range("N100").Value = "$/docena de paquetes"
$ws.Range("O100").Value = "Provincia de Cautín"
$ws.Range("P100").Value = 667
$ws.Range("Q100").Value = 12
$ws.Range("R100").Value = "Hortaliza"
